$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.301.34'
$ws.Range("E2").Value = '  +1.63%  '

$ws.Range("D3").Value = '3.938.58'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '492.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.13%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.732'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("E10").Value = '  +3.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000351'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.32'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.21%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.79%  '

$ws.Range("D14").Value = '4.568.95'

$ws.Range("D15").Value = '3.951.39'
$ws.Range("E15").Value = '  +1.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.17%  '

$ws.Range("E17").Value = '  -0.86%  '

$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.93%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.60%  '

$ws.Range("D20").Value = '69.380.71'
$ws.Range("E20").Value = '  +1.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '440.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.14%  '

$ws.Range("E26").Value = '  +4.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '705.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.86%  '

$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.44%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.130'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.37%  '

$ws.Range("E33").Value = '  +0.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.464'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +17.43%  '

$ws.Range("D35").Value = '0.0₃0912'
$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.09'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '61.65'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.03%  '

$ws.Range("E39").Value = '  +1.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("E43").Value = '  -1.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.18%  '

$ws.Range("E47").Value = '  +7.56%  '

$ws.Range("D48").Value = '0.0₆0357'
$ws.Range("E48").Value = '  +1.36%  '

$ws.Range("E49").Value = '  +5.51%  '

$ws.Range("E50").Value = '  -0.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.72%  '
